# Plantilla de Casos de Uso - actualizacion de esfuerzos (rubro de Grupos)
# y de la vista activa de la hoja, segun el commit:
# "Diagramas de secuencia, actualizacion plantillas"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Esfuerzo (hrs) actualizado para los CU de Grupos (CREAR/EDITAR/ELIMINAR/CONSULTAR GRUPOS)
$ws.Range("F14").Value = 1.5
$ws.Range("F15").Value = 0.83
$ws.Range("F16").Value = 1.086
$ws.Range("F17").Value = 1.65

# La hoja vuelve a quedar posicionada/seleccionada sobre la fila de los CU de Grupos
$ws.Activate() | Out-Null
$ws.Range("F17").Select() | Out-Null
